# Add working set of sequences
# For every row whose "cued" column (E) is TRUE, the row previously only had
# columns A-F populated. This fills in columns G-N (image2..image9 style
# columns) with the same "N/A" placeholder value already used in column F,
# mirroring the fully-populated rows (where E is FALSE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 7, 10, 13, 14, 17, 19, 23, 25, 30, 35, 37, 38, 44, 45, 51, 53, 55, 56, 59, 65, 70, 71, 78, 79, 85, 87, 91, 94, 97, 98, 103, 108, 110, 113, 117, 119, 120, 122, 123, 125, 126, 129, 135, 139, 142, 144, 149, 155, 156, 157, 159, 160, 162, 163, 168, 170, 172, 178, 181, 183, 186, 187, 189)

foreach ($r in $rows) {
    $fill = $ws.Cells.Item($r, 6).Value2
    for ($c = 7; $c -le 14; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $fill
    }
}
